# Update the NATMI Wnt5b-Fzd8 LR-pair sheet with newly recomputed TPM-based
# statistics, and add the missing "MuSCs" sending-cluster rows.
#
# The table has columns:
#   A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
#   E..T numeric NATMI metrics
#
# Rows 2-4 (existing, Sending cluster = FAPs) get refreshed numeric values,
# and three new rows (5-7, Sending cluster = MuSCs) are appended.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row data: Sending, Ligand, Receptor, Target, then E..T numeric metrics
# ---------------------------------------------------------------------
$rows = @(
    @{ Row=2; A="FAPs";  D="ECs";   E=1; F=0.3333333333333333; G=0.4046233333333333; H=1.21387;
       I=0.4348530100317075; J=0.4348530100317076; K=3; L=1;
       M=3.083576666666666;  N=9.250729999999999;
       O=0.2272509363535097; P=0.2272509363535097;
       Q=1.247687069455555;  R=11.2291836251;
       S=0.09882075370584766; T=0.09882075370584767 },

    @{ Row=3; A="FAPs";  D="FAPs";  E=1; F=0.3333333333333333; G=0.4046233333333333; H=1.21387;
       I=0.4348530100317075; J=0.4348530100317076; K=3; L=1;
       M=6.453984666666667;  N=19.361954;
       O=0.4756405360586227; P=0.4756405360586227;
       Q=2.611432789108889;  R=23.50289510198;
       S=0.206833718798187;  T=0.206833718798187 },

    @{ Row=4; A="FAPs";  D="MuSCs"; E=1; F=0.3333333333333333; G=0.4046233333333333; H=1.21387;
       I=0.4348530100317075; J=0.4348530100317076; K=3; L=1;
       M=4.031477000000001;  N=12.094431;
       O=0.2971085275878677; P=0.2971085275878677;
       Q=1.631229661996667;  R=14.68106695797;
       S=0.1291985375276729; T=0.1291985375276729 },

    @{ Row=5; A="MuSCs"; D="ECs";   E=2; F=0.6666666666666666; G=0.5258596666666667; H=1.577579;
       I=0.5651469899682925; J=0.5651469899682925; K=3; L=1;
       M=3.083576666666666;  N=9.250729999999999;
       O=0.2272509363535097; P=0.2272509363535097;
       Q=1.621528598074444;  R=14.59375738267;
       S=0.128430182647662;  T=0.128430182647662 },

    @{ Row=6; A="MuSCs"; D="FAPs";  E=2; F=0.6666666666666666; G=0.5258596666666667; H=1.577579;
       I=0.5651469899682925; J=0.5651469899682925; K=3; L=1;
       M=6.453984666666667;  N=19.361954;
       O=0.4756405360586227; P=0.4756405360586227;
       Q=3.393890225485112;  R=30.545012029366;
       S=0.2688068172604357; T=0.2688068172604357 },

    @{ Row=7; A="MuSCs"; D="MuSCs"; E=2; F=0.6666666666666666; G=0.5258596666666667; H=1.577579;
       I=0.5651469899682925; J=0.5651469899682925; K=3; L=1;
       M=4.031477000000001;  N=12.094431;
       O=0.2971085275878677; P=0.2971085275878677;
       Q=2.119991151394334;  R=19.079920362549;
       S=0.1679099900601948; T=0.1679099900601948 }
)

foreach ($r in $rows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value  = $r.A       # Sending cluster
    $ws.Cells.Item($row, 2).Value  = "Wnt5b"    # Ligand symbol
    $ws.Cells.Item($row, 3).Value  = "Fzd8"     # Receptor symbol
    $ws.Cells.Item($row, 4).Value  = $r.D       # Target cluster

    $ws.Cells.Item($row, 5).Value  = $r.E
    $ws.Cells.Item($row, 6).Value  = $r.F
    $ws.Cells.Item($row, 7).Value  = $r.G
    $ws.Cells.Item($row, 8).Value  = $r.H
    $ws.Cells.Item($row, 9).Value  = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = $r.R
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $r.T
}

Write-Output "Updated rows 2-7 with refreshed TPM-based NATMI statistics"
